$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1482.9166
$ws.Range("I43").Value = 900
$ws.Range("K43").Value = 900
$ws.Range("M43").Value = -831
$ws.Range("H96").Value = 2153
$ws.Range("I96").Value = 2153
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 6459
$ws.Range("L96").Value = 0
$ws.Range("M96").ClearContents()
$ws.Range("N96").Value = -5086
$ws.Range("H98").Value = 2958.9167
$ws.Range("I98").Value = 3137
$ws.Range("K98").Value = 3137
$ws.Range("M98").Value = -1639
$ws.Range("H116").Value = 12883.538
$ws.Range("I116").Value = 35133.332
$ws.Range("J116").Value = 6208.6
$ws.Range("K116").Value = 35133.332
$ws.Range("L116").Value = 6208.6
$ws.Range("M116").Value = -31691.332
$ws.Range("N116").Value = -13092.6
$ws.Range("H122").Value = 2958.9167
$ws.Range("I122").Value = 3137
$ws.Range("K122").Value = 9411
$ws.Range("M122").Value = -6961
$ws.Range("H137").Value = 1438.5758
$ws.Range("J137").Value = 2091.4666
$ws.Range("L137").Value = 6274.399800000001
$ws.Range("N137").Value = -11374.3998
$ws.Range("H138").Value = 2401.4102
$ws.Range("I138").Value = 2627
$ws.Range("J138").Value = 2040.4667
$ws.Range("K138").Value = 7881
$ws.Range("L138").Value = 6121.4001
$ws.Range("M138").Value = -2741
$ws.Range("N138").Value = -16401.4001
$ws.Range("H139").Value = 47685.715
$ws.Range("J139").Value = 47685.715
$ws.Range("L139").Value = 47685.715
$ws.Range("N139").Value = -57965.715
$ws.Range("H140").Value = 53151.867
$ws.Range("J140").Value = 53151.867
$ws.Range("L140").Value = 53151.867
$ws.Range("N140").Value = -63511.867
$ws.Range("H141").Value = 4995.6665
$ws.Range("I141").Value = 3995
$ws.Range("K141").Value = 11985
$ws.Range("M141").Value = -6805
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3877.7415
$ws.Range("I32").Value = 2471.976
$ws.Range("J32").Value = 7567.875
$ws.Range("K32").Value = 2471.976
$ws.Range("L32").Value = 7567.875
$ws.Range("M32").Value = -2184.976
$ws.Range("N32").Value = -8141.875
$ws.Range("H45").Value = 1481.8182
$ws.Range("I45").Value = 989.2
$ws.Range("K45").Value = 989.2
$ws.Range("M45").Value = -612.2
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1919.5
$ws.Range("I107").Value = 1799.4445
$ws.Range("J107").Value = 3000
$ws.Range("K107").Value = 1799.4445
$ws.Range("L107").Value = 3000
$ws.Range("M107").Value = 120.5554999999999
$ws.Range("N107").Value = -6840
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2712.1428
$ws.Range("I31").Value = 1328.4375
$ws.Range("J31").Value = 7140
$ws.Range("K31").Value = 1328.4375
$ws.Range("L31").Value = 7140
$ws.Range("M31").Value = -1033.4375
$ws.Range("N31").Value = -7730
$ws.Range("H34").Value = 2712.1428
$ws.Range("I34").Value = 1328.4375
$ws.Range("J34").Value = 7140
$ws.Range("K34").Value = 1328.4375
$ws.Range("L34").Value = 7140
$ws.Range("M34").Value = -1126.4375
$ws.Range("N34").Value = -7544
$ws.Range("H43").Value = 38998
$ws.Range("J43").Value = 38998
$ws.Range("L43").Value = 38998
$ws.Range("N43").Value = -39366
$ws.Range("H101").Value = 38998
$ws.Range("J101").Value = 38998
$ws.Range("L101").Value = 38998
$ws.Range("N101").Value = -45488
$ws.Range("H132").Value = 2507.7856
$ws.Range("I132").Value = 1699.5714
$ws.Range("J132").Value = 3316
$ws.Range("K132").Value = 5098.7142
$ws.Range("L132").Value = 9948
$ws.Range("M132").Value = -2568.7142
$ws.Range("N132").Value = -15008
$ws.Range("H134").Value = 1138.2
$ws.Range("I134").Value = 922.75
$ws.Range("J134").Value = 2000
$ws.Range("K134").Value = 2768.25
$ws.Range("L134").Value = 6000
$ws.Range("M134").Value = -233.25
$ws.Range("N134").Value = -11070
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 769.8125
$ws.Range("J107").Value = 800.93335
$ws.Range("L107").Value = 2402.80005
$ws.Range("N107").Value = -6242.80005
$ws.Range("H131").Value = 773.83
$ws.Range("J131").Value = 810.23914
$ws.Range("L131").Value = 2430.71742
$ws.Range("N131").Value = -12510.71742
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2628.9473
$ws.Range("I102").Value = 3676.1667
$ws.Range("J102").Value = 2145.6155
$ws.Range("K102").Value = 3676.1667
$ws.Range("L102").Value = 2145.6155
$ws.Range("M102").Value = -2054.1667
$ws.Range("N102").Value = -5389.6155
$ws.Range("H122").Value = 2433.6155
$ws.Range("J122").Value = 2693
$ws.Range("L122").Value = 8079
$ws.Range("N122").Value = -12979
$ws.Range("H126").Value = 65405.375
$ws.Range("I126").Value = 3389
$ws.Range("J126").Value = 127421.75
$ws.Range("K126").Value = 10167
$ws.Range("L126").Value = 382265.25
$ws.Range("M126").Value = -7697
$ws.Range("N126").Value = -387205.25
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1002
$ws.Range("I22").Value = 1002
$ws.Range("K22").Value = 1002
$ws.Range("M22").Value = -707
$ws.Range("H27").Value = 1002
$ws.Range("I27").Value = 1002
$ws.Range("K27").Value = 1002
$ws.Range("M27").Value = -895
$ws.Range("H40").Value = 13442.111
$ws.Range("I40").Value = 11992.5
$ws.Range("J40").Value = 13856.286
$ws.Range("K40").Value = 11992.5
$ws.Range("L40").Value = 13856.286
$ws.Range("M40").Value = -11856.5
$ws.Range("N40").Value = -14128.286
$ws.Range("H122").Value = 6950.5
$ws.Range("J122").Value = 8000.8335
$ws.Range("L122").Value = 24002.5005
$ws.Range("N122").Value = -28902.5005
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 30000
$ws.Range("J70").Value = 30000
$ws.Range("L70").Value = 30000
$ws.Range("N70").Value = -30630
$ws.Range("H73").Value = 30000
$ws.Range("J73").Value = 30000
$ws.Range("L73").Value = 30000
$ws.Range("N73").Value = -32184
$ws.Range("H95").Value = 49999
$ws.Range("J95").Value = 49999
$ws.Range("L95").Value = 49999
$ws.Range("N95").Value = -55491
$ws.Range("H122").Value = 184146.72
$ws.Range("I122").Value = 320506.75
$ws.Range("J122").Value = 2333.3333
$ws.Range("K122").Value = 961520.25
$ws.Range("L122").Value = 6999.999899999999
$ws.Range("M122").Value = -959070.25
$ws.Range("N122").Value = -11899.9999
$ws.Range("H132").Value = 8673.956
$ws.Range("I132").Value = 1504
$ws.Range("K132").Value = 4512
$ws.Range("M132").Value = -1982
